$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 34.32326441073354
$ws.Range("G2").Value = 33.93231324440875
$ws.Range("H2").Value = 34.71324539606103
$ws.Range("I2").Value = 3.833116648445406
$ws.Range("J2").Value = 3.805656843724844
$ws.Range("K2").Value = 3.86061106235283
$ws.Range("L2").Value = 0.2688364430516824
$ws.Range("M2").Value = 0.2667919453456743
$ws.Range("N2").Value = 0.2708875651108002
$ws.Range("F3").Value = 0.0002086230621730263
$ws.Range("G3").Value = 0.00000001256795025251431
$ws.Range("H3").Value = 0.0005862637385339242
$ws.Range("I3").Value = 0.0001958438724335321
$ws.Range("J3").Value = 0.00000001172202137309654
$ws.Range("K3").Value = 0.0005503416370005344
$ws.Range("L3").Value = 0.000209302570309247
$ws.Range("M3").Value = 0.00000001254978277096446
$ws.Range("N3").Value = 0.0005883541305949057
$ws.Range("F4").Value = 34.32347303379571
$ws.Range("G4").Value = 33.9323132569767
$ws.Range("H4").Value = 34.71383165979957
$ws.Range("I4").Value = 3.83331249231784
$ws.Range("J4").Value = 3.805656855446865
$ws.Range("K4").Value = 3.86116140398983
$ws.Range("L4").Value = 0.2690457456219916
$ws.Range("M4").Value = 0.2667919578954571
$ws.Range("N4").Value = 0.2714759192413951
